# Weekly update: insert a new price-report row for the week right after the
# existing row 44 (D44 = 2021-06-22 ... stays put), pushing the former rows
# 45-60 down by one (to 46-61), and populate the newly inserted row 45 with
# the new week's data. All other columns for this market/category are
# constant across the block, so they are copied through unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 45 - shifts existing rows 45:60 down to 46:61
# and copies formatting (incl. the date number format) from the row above.
$ws.Rows.Item(45).Insert()

$ws.Range("A45").Value = 8
$ws.Range("B45").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C45").Value = 'Coquimbo'
$ws.Range("D45").Value = 44523
$ws.Range("E45").Value = 4
$ws.Range("F45").Value = 100112052
$ws.Range("G45").Value = 'Albahaca'
$ws.Range("H45").Value = 'Sin especificar'
$ws.Range("I45").Value = 'Primera'
$ws.Range("J45").Value = 760
$ws.Range("K45").Value = 3000
$ws.Range("L45").Value = 4000
$ws.Range("M45").Value = 3500
$ws.Range("N45").Value = '$/paquete'
$ws.Range("O45").Value = 'Región de Arica y Parinacota'
$ws.Range("P45").Value = 3500
$ws.Range("Q45").Value = 1
$ws.Range("R45").Value = 'Hortaliza'
